# Applies the cryptos-list price/volume refresh described in the commit diff.
# Numeric-looking values are written with a leading quote-prefix (Set-CellText)
# so Excel keeps them as text (matching the source inlineStr cells) instead of
# silently parsing them into numbers and losing formatting (e.g. '1.00' -> 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($Sheet, [string]$Address, [string]$Text, [bool]$ForceText)
    if ($ForceText) {
        $Sheet.Range($Address).Value = "'" + $Text
    } else {
        $Sheet.Range($Address).Value = $Text
    }
}

Set-CellText $ws 'D2' '34.582.11' $false
Set-CellText $ws 'E2' '  +1.01%  ' $false
Set-CellText $ws 'D3' '1.818.71' $false
Set-CellText $ws 'E3' '  +1.63%  ' $false
Set-CellText $ws 'D5' '228.29' $true
Set-CellText $ws 'E5' '  +0.94%  ' $false
Set-CellText $ws 'E6' '  +1.31%  ' $false
Set-CellText $ws 'E7' '  +0.02%  ' $false
Set-CellText $ws 'D8' '34.99' $true
Set-CellText $ws 'E8' '  +8.49%  ' $false
Set-CellText $ws 'E9' '  +1.77%  ' $false
Set-CellText $ws 'E10' '  +0.96%  ' $false
Set-CellText $ws 'D11' '0.0950' $true
Set-CellText $ws 'E11' '  +0.36%  ' $false
Set-CellText $ws 'D12' '2.082.39' $false
Set-CellText $ws 'E12' '  +1.67%  ' $false
Set-CellText $ws 'D13' '11.47' $true
Set-CellText $ws 'E13' '  +3.77%  ' $false
Set-CellText $ws 'D14' '1.835.83' $false
Set-CellText $ws 'E14' '  +2.65%  ' $false
Set-CellText $ws 'E15' '  +3.32%  ' $false
Set-CellText $ws 'D16' '34.581.41' $false
Set-CellText $ws 'E16' '  +1.05%  ' $false
Set-CellText $ws 'D17' '4.34' $true
Set-CellText $ws 'E17' '  +3.53%  ' $false
Set-CellText $ws 'D18' '69.08' $true
Set-CellText $ws 'E18' '  +1.60%  ' $false
Set-CellText $ws 'B19' 'ShibaInu' $false
Set-CellText $ws 'C19' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib' $false
Set-CellText $ws 'D19' '0.0₃0802' $false
Set-CellText $ws 'E19' '  +0.40%  ' $false
Set-CellText $ws 'B20' 'BitcoinCash' $false
Set-CellText $ws 'C20' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch' $false
Set-CellText $ws 'D20' '246.86' $true
Set-CellText $ws 'E20' '  +0.38%  ' $false
Set-CellText $ws 'D21' '11.52' $true
Set-CellText $ws 'E21' '  +5.28%  ' $false
Set-CellText $ws 'D22' '1.00' $true
Set-CellText $ws 'E22' '  +0.09%  ' $false
Set-CellText $ws 'D23' '4.20' $true
Set-CellText $ws 'E23' '  +1.16%  ' $false
Set-CellText $ws 'D24' '171.37' $true
Set-CellText $ws 'E24' '  +5.98%  ' $false
Set-CellText $ws 'D25' '2.08' $true
Set-CellText $ws 'E25' '  +1.19%  ' $false
Set-CellText $ws 'D26' '7.40' $true
Set-CellText $ws 'E26' '  +3.05%  ' $false
Set-CellText $ws 'D27' '16.77' $true
Set-CellText $ws 'E27' '  +2.66%  ' $false
Set-CellText $ws 'D28' '0.116' $true
Set-CellText $ws 'E28' '  +1.45%  ' $false
Set-CellText $ws 'E29' '  -0.43%  ' $false
Set-CellText $ws 'E30' '  +7.22%  ' $false
Set-CellText $ws 'D31' '0.0532' $true
Set-CellText $ws 'E31' '  +2.10%  ' $false
Set-CellText $ws 'E32' '  +2.68%  ' $false
Set-CellText $ws 'D33' '1.25' $true
Set-CellText $ws 'E33' '  +0.97%  ' $false
Set-CellText $ws 'E34' '  +2.57%  ' $false
Set-CellText $ws 'E35' '  +1.83%  ' $false
Set-CellText $ws 'D36' '1.421.23' $false
Set-CellText $ws 'E36' '  -1.29%  ' $false
Set-CellText $ws 'E37' '  +2.44%  ' $false
Set-CellText $ws 'D38' '1.06' $true
Set-CellText $ws 'E38' '  +0.95%  ' $false
Set-CellText $ws 'B39' 'VeChain' $false
Set-CellText $ws 'C39' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' $false
Set-CellText $ws 'D39' '0.0191' $true
Set-CellText $ws 'E39' '  +0.64%  ' $false
Set-CellText $ws 'B40' 'Aave' $false
Set-CellText $ws 'C40' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave' $false
Set-CellText $ws 'D40' '86.00' $true
Set-CellText $ws 'E40' '  +4.94%  ' $false
Set-CellText $ws 'D41' '2.85' $true
Set-CellText $ws 'E41' '  +4.37%  ' $false
Set-CellText $ws 'E42' '  +3.80%  ' $false
Set-CellText $ws 'E43' '  +1.12%  ' $false
Set-CellText $ws 'D44' '13.96' $true
Set-CellText $ws 'E44' '  -0.89%  ' $false
Set-CellText $ws 'D45' '0.0528' $true
Set-CellText $ws 'E45' '  +1.65%  ' $false
Set-CellText $ws 'E46' '  +2.98%  ' $false
Set-CellText $ws 'E47' '  +0.19%  ' $false
Set-CellText $ws 'D48' '1.982.84' $false
Set-CellText $ws 'E48' '  +1.97%  ' $false
Set-CellText $ws 'D49' '105.89' $true
Set-CellText $ws 'E49' '  +0.48%  ' $false
Set-CellText $ws 'D50' '0.0₆0132' $false
Set-CellText $ws 'E50' '  +1.98%  ' $false
Set-CellText $ws 'E51' '  -0.04%  ' $false
